$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.112737948256028
$ws.Range("C2").Value = 0.426895895633562
$ws.Range("D2").Value = 0.1572417756599283
$ws.Range("E2").Value = 0.7087734795590415
$ws.Range("F2").Value = 0.482846718646189
$ws.Range("G2").Value = 0.5906624340427387
$ws.Range("H2").Value = 0.4942628133784369
$ws.Range("I2").Value = 0.9976736988314454
$ws.Range("J2").Value = 1.256195640754894
$ws.Range("K2").Value = -0.4677869442234193

$ws.Range("B3").Value = 0.9584078752945684
$ws.Range("C3").Value = 1.509939579193682
$ws.Range("D3").Value = 1.284012818280829
$ws.Range("E3").Value = 1.391828533677379
$ws.Range("F3").Value = 1.295428913013077
$ws.Range("G3").Value = 1.798839798466086
$ws.Range("H3").Value = 2.057361740389535
$ws.Range("I3").Value = 0.3333791554112209
$ws.Range("J3").Value = 0.9870183717739222
$ws.Range("K3").Value = 1.350562306682353

$ws.Range("B4").Value = 0.2878486627132291
$ws.Range("C4").Value = 0.3956643781097787
$ws.Range("D4").Value = 0.299264757445477
$ws.Range("E4").Value = 0.8026756428984855
$ws.Range("F4").Value = 1.061197584821935
$ws.Range("G4").Value = -0.6627850001563792
$ws.Range("H4").Value = -0.00914578379367792
$ws.Range("I4").Value = 0.3543981511147533
$ws.Range("J4").Value = -0.5463776937648762
$ws.Range("K4").Value = 0.5239978631887766

$ws.Range("B5").Value = 0.3064042510528396
$ws.Range("C5").Value = 0.8098151365058481
$ws.Range("D5").Value = 1.068337078429297
$ws.Range("E5").Value = -0.6556455065490165
$ws.Range("F5").Value = -0.002006290186315263
$ws.Range("G5").Value = 0.3615376447221159
$ws.Range("H5").Value = -0.5392382001575136
$ws.Range("I5").Value = 0.5311373567961393
$ws.Range("J5").Value = -0.04592255722032257
$ws.Range("K5").Value = -0.2515285413133388

$ws.Range("B6").Value = 0.7921607576674977
$ws.Range("C6").Value = -0.9318218273108161
$ws.Range("D6").Value = -0.2781826109481148
$ws.Range("E6").Value = 0.0853613239603164
$ws.Range("F6").Value = -0.8154145209193131
$ws.Range("G6").Value = 0.2549610360343397
$ws.Range("H6").Value = -0.3220988779821221
$ws.Range("I6").Value = -0.5277048620751383
$ws.Range("J6").Value = 0.4514401196459589
$ws.Range("K6").Value = -0.2849062236914031

$ws.Range("B7").Value = -0.4488784752138651
$ws.Range("C7").Value = -0.08533454030543387
$ws.Range("D7").Value = -0.9861103851850633
$ws.Range("E7").Value = 0.08426517176858947
$ws.Range("F7").Value = -0.4927947422478724
$ws.Range("G7").Value = -0.6984007263408886
$ws.Range("H7").Value = 0.2807442553802086
$ws.Range("I7").Value = -0.4556020879571533
$ws.Range("J7").Value = -0.4735719026113467
$ws.Range("K7").Value = -0.4501670370710099

$ws.Range("B8").Value = -0.7270750933112256
$ws.Range("C8").Value = 0.3433004636424272
$ws.Range("D8").Value = -0.2337594503740346
$ws.Range("E8").Value = -0.4393654344670508
$ws.Range("F8").Value = 0.5397795472540464
$ws.Range("G8").Value = -0.1965667960833156
$ws.Range("H8").Value = -0.2145366107375089
$ws.Range("I8").Value = -0.1911317451971721
$ws.Range("J8").Value = -0.5563061809639129
$ws.Range("K8").Value = 0.09951333620703928

$ws.Range("B9").Value = 0.2618430268646463
$ws.Range("C9").Value = 0.05623704277163011
$ws.Range("D9").Value = 1.035382024492727
$ws.Range("E9").Value = 0.2990356811553654
$ws.Range("F9").Value = 0.281065866501172
$ws.Range("G9").Value = 0.3044707320415089
$ws.Range("H9").Value = -0.06070370372523193
$ws.Range("I9").Value = 0.5951158134457202
$ws.Range("J9").Value = 0.4642449325042965
$ws.Range("K9").Value = 0.3547271335801747

$ws.Range("B10").Value = 0.2836965368615766
$ws.Range("C10").Value = -0.4526498064757853
$ws.Range("D10").Value = -0.4706196211299787
$ws.Range("E10").Value = -0.4472147555896419
$ws.Range("F10").Value = -0.8123891913563827
$ws.Range("G10").Value = -0.1565696741854305
$ws.Range("H10").Value = -0.2874405551268542
$ws.Range("I10").Value = -0.3969583540509761
$ws.Range("J10").Value = -1.226736141118529
$ws.Range("K10").Value = -0.6462103519498816

$ws.Range("B11").Value = -0.2645568503594945
$ws.Range("C11").Value = -0.2411519848191577
$ws.Range("D11").Value = -0.6063264205858985
$ws.Range("E11").Value = 0.04949309658505369
$ws.Range("F11").Value = -0.08137778435637
$ws.Range("G11").Value = -0.1908955832804919
$ws.Range("H11").Value = -1.020673370348045
$ws.Range("I11").Value = -0.4401475811793975
$ws.Range("J11").Value = -0.5235336826091774
$ws.Range("K11").Value = -0.7414435601489361

$ws.Range("B12").Value = -0.5627004823385774
$ws.Range("C12").Value = 0.09311903483237477
$ws.Range("D12").Value = -0.03775184610904891
$ws.Range("E12").Value = -0.1472696450331708
$ws.Range("F12").Value = -0.9770474321007241
$ws.Range("G12").Value = -0.3965216429320764
$ws.Range("H12").Value = -0.4799077443618563
$ws.Range("I12").Value = -0.697817621901615
$ws.Range("J12").Value = 0.8659601817345554
$ws.Range("K12").Value = -0.4700235243475083

$ws.Range("B13").Value = 0.3116331488582624
$ws.Range("C13").Value = 0.2021153499341405
$ws.Range("D13").Value = -0.6276624371334127
$ws.Range("E13").Value = -0.04713664796476502
$ws.Range("F13").Value = -0.130522749394545
$ws.Range("G13").Value = -0.3484326269343037
$ws.Range("H13").Value = 1.215345176701867
$ws.Range("I13").Value = -0.1206385293801969
$ws.Range("J13").Value = 1.319975887149931
$ws.Range("K13").Value = 0.2195317728891086

$ws.Range("B14").Value = -1.101838862130705
$ws.Range("C14").Value = -0.5213130729620578
$ws.Range("D14").Value = -0.6046991743918377
$ws.Range("E14").Value = -0.8226090519315964
$ws.Range("F14").Value = 0.741168751704574
$ws.Range("G14").Value = -0.5948149543774897
$ws.Range("H14").Value = 0.8457994621526386
$ws.Range("I14").Value = -0.2546446521081841
$ws.Range("J14").Value = -0.464919050277854
$ws.Range("K14").Value = 0.1336251531058764

$ws.Range("B15").Value = -0.1899379683343848
$ws.Range("C15").Value = -0.4078478458741435
$ws.Range("D15").Value = 1.155929957762027
$ws.Range("E15").Value = -0.1800537483200367
$ws.Range("F15").Value = 1.260560668210092
$ws.Range("G15").Value = 0.1601165539492688
$ws.Range("H15").Value = -0.05015784422040104
$ws.Range("I15").Value = 0.5483863591633293
$ws.Range("J15").Value = 0.2227177183881149
$ws.Range("K15").Value = 0.9398975197637186

$ws.Range("B16").Value = 1.374209598179478
$ws.Range("C16").Value = 0.03822589209741434
$ws.Range("D16").Value = 1.478840308627543
$ws.Range("E16").Value = 0.3783961943667199
$ws.Range("F16").Value = 0.16812179619705
$ws.Range("G16").Value = 0.7666659995807804
$ws.Range("H16").Value = 0.440997358805566
$ws.Range("I16").Value = 1.15817716018117
$ws.Range("J16").Value = 2.854049435818871
$ws.Range("K16").Value = 9.698296230743448

$ws.Range("B17").Value = -0.1789830073774904
$ws.Range("C17").Value = 1.261631409152638
$ws.Range("D17").Value = 0.1611872948918152
$ws.Range("E17").Value = -0.0490871032778547
$ws.Range("F17").Value = 0.5494571001058757
$ws.Range("G17").Value = 0.2237884593306613
$ws.Range("H17").Value = 0.9409682607062649
$ws.Range("I17").Value = 2.636840536343966
$ws.Range("J17").Value = 9.481087331268544
$ws.Range("K17").Value = -8.069566751832767

$ws.Range("B18").Value = 0.9519452398945764
$ws.Range("C18").Value = -0.1484988743662463
$ws.Range("D18").Value = -0.3587732725359162
$ws.Range("E18").Value = 0.2397709308478142
$ws.Range("F18").Value = -0.0858977099274002
$ws.Range("G18").Value = 0.6312820914482035
$ws.Range("H18").Value = 2.327154367085904
$ws.Range("I18").Value = 9.171401162010483
$ws.Range("J18").Value = -8.379252921090828
$ws.Range("K18").Value = -0.6698563214346814

$ws.Range("B19").Value = 0.005805259465757717
$ws.Range("C19").Value = -0.2044691387039121
$ws.Range("D19").Value = 0.3940750646798182
$ws.Range("E19").Value = 0.06840642390460383
$ws.Range("F19").Value = 0.7855862252802075
$ws.Range("G19").Value = 2.481458500917908
$ws.Range("H19").Value = 9.325705295842486
$ws.Range("I19").Value = -8.224948787258825
$ws.Range("J19").Value = -0.5155521876026774
$ws.Range("K19").Value = 1.067360473284234

$ws.Range("B20").Value = -0.8548955599230954
$ws.Range("C20").Value = -0.256351356539365
$ws.Range("D20").Value = -0.5820199973145794
$ws.Range("E20").Value = 0.1351598040610243
$ws.Range("F20").Value = 1.831032079698725
$ws.Range("G20").Value = 8.675278874623302
$ws.Range("H20").Value = -8.875375208478008
$ws.Range("I20").Value = -1.165978608821861
$ws.Range("J20").Value = 0.416934052065051
$ws.Range("K20").Value = -2.549352940284205

$ws.Range("B21").Value = -0.04018630792115581
$ws.Range("C21").Value = -0.3658549486963703
$ws.Range("D21").Value = 0.3513248526792334
$ws.Range("E21").Value = 2.047197128316934
$ws.Range("F21").Value = 8.891443923241512
$ws.Range("G21").Value = -8.659210159859798
$ws.Range("H21").Value = -0.9498135602036515
$ws.Range("I21").Value = 0.6330991006832601
$ws.Range("J21").Value = -2.333187891665996
$ws.Range("K21").Value = -0.06915673700153271

$ws.Range("B22").Value = 0.04117418033858611
$ws.Range("C22").Value = 0.7583539817141898
$ws.Range("D22").Value = 2.454226257351891
$ws.Range("E22").Value = 9.298473052276469
$ws.Range("F22").Value = -8.252181030824842
$ws.Range("G22").Value = -0.5427844311686951
$ws.Range("H22").Value = 1.040128229718217
$ws.Range("I22").Value = -1.926158762631039
$ws.Range("J22").Value = 0.3378723920334236
$ws.Range("K22").Value = -0.113281762893315

$ws.Range("B23").Value = 0.25931176453677
$ws.Range("C23").Value = 1.955184040174471
$ws.Range("D23").Value = 8.799430835099049
$ws.Range("E23").Value = -8.751223248002262
$ws.Range("F23").Value = -1.041826648346115
$ws.Range("G23").Value = 0.5410860125407967
$ws.Range("H23").Value = -2.425200979808459
$ws.Range("I23").Value = -0.1611698251439962
$ws.Range("J23").Value = -0.6123239800707349
$ws.Range("K23").Value = -0.752688043584709

$ws.Range("B24").Value = 2.114032953329983
$ws.Range("C24").Value = 8.958279748254562
$ws.Range("D24").Value = -8.592374334846749
$ws.Range("E24").Value = -0.8829777351906024
$ws.Range("F24").Value = 0.6999349256963092
$ws.Range("G24").Value = -2.266352066652947
$ws.Range("H24").Value = -0.002320911988483623
$ws.Range("I24").Value = -0.4534750669152223
$ws.Range("J24").Value = -0.5938391304291964
$ws.Range("K24").Value = -0.06760011956697885
